$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 689, shifting existing rows 689:739 down to 690:740
$ws.Rows(689).Insert()

# Populate the newly inserted row 689 with the new record's data
$ws.Cells.Item(689, 1).Value = 10
$ws.Cells.Item(689, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(689, 3).Value = "La Araucanía"
$ws.Cells.Item(689, 4).Value = 45166
$ws.Cells.Item(689, 5).Value = 9
$ws.Cells.Item(689, 6).Value = 100112028
$ws.Cells.Item(689, 7).Value = "Sandia"
$ws.Cells.Item(689, 8).Value = "Sin especificar"
$ws.Cells.Item(689, 9).Value = "Primera"
$ws.Cells.Item(689, 10).Value = 700
$ws.Cells.Item(689, 11).Value = 1300
$ws.Cells.Item(689, 12).Value = 1400
$ws.Cells.Item(689, 13).Value = 1357
$ws.Cells.Item(689, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(689, 15).Value = "Perú"
$ws.Cells.Item(689, 16).Value = 1357
$ws.Cells.Item(689, 17).Value = 1
$ws.Cells.Item(689, 18).Value = "Hortaliza"
